$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> updated Price (D) and Volume(1h) (E) values.
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the original inlineStr/text cell contents)
# instead of auto-converting numeric-looking strings or percentages.
$updates = @{
    2 = @{ D="321.93"; E="5.20%" }
    3 = @{ D="36.20"; E="-0.03%" }
    4 = @{ D="5.123"; E="1.60%" }
    5 = @{ D="0.08071"; E="3.24%" }
    6 = @{ D="2.153"; E="1.33%" }
    7 = @{ D="8.067"; E="1.97%" }
    8 = @{ E="1.31%" }
    9 = @{ D="0.9295"; E="1.13%" }
    10 = @{ D="0.1011"; E="5.46%" }
    11 = @{ D="0.1883"; E="0.72%" }
    12 = @{ D="0.09249"; E="6.61%" }
    13 = @{ D="0.03571"; E="1.74%" }
    14 = @{ D="0.09934"; E="0.06%" }
    15 = @{ E="0.56%" }
    16 = @{ D="0.005665"; E="-0.52%" }
    17 = @{ D="3.458"; E="0.04%" }
    18 = @{ E="16.70%" }
    19 = @{ E="-1.33%" }
    20 = @{ D="0.1330"; E="-1.40%" }
    21 = @{ D="5.064"; E="6.38%" }
    22 = @{ E="-4.07%" }
    23 = @{ D="0.04605"; E="-0.18%" }
    24 = @{ E="0.95%" }
    25 = @{ E="-6.97%" }
    26 = @{ D="0.0001300"; E="-7.13%" }
    27 = @{ D="0.0004501"; E="65.12%" }
    39 = @{ D="0.01962"; E="8.12%" }
    40 = @{ D="0.04980"; E="4.89%" }
    41 = @{ D="0.007798"; E="4.26%" }
    42 = @{ D="0.1401"; E="0.11%" }
    43 = @{ D="0.007692"; E="-0.70%" }
    44 = @{ D="0.002067"; E="-7.70%" }
    45 = @{ D="0.01181"; E="9.41%" }
    46 = @{ D="0.00006385"; E="2.22%" }
    47 = @{ D="0.00000000750"; E="-0.29%" }
    48 = @{ E="-11.93%" }
    49 = @{ D="0.001901"; E="-5.27%" }
    50 = @{ E="-0.29%" }
    51 = @{ E="-0.29%" }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("D")) {
        $ws.Range("D$row").Value = "'" + $rowData.D
    }
    if ($rowData.ContainsKey("E")) {
        $ws.Range("E$row").Value = "'" + $rowData.E
    }
}
